$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.196.49"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "2.175.65"
$ws.Range("E3").Value = "  -0.46%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "254.48"
$ws.Range("D5").Style = $ws.Range("B2").Style
$ws.Range("E5").Value = "  +5.34%  "
$ws.Range("E6").Value = "  +1.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "67.65"
$ws.Range("D7").Style = $ws.Range("B2").Style
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.572"
$ws.Range("D9").Style = $ws.Range("B2").Style
$ws.Range("E9").Value = "  +6.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.44"
$ws.Range("D10").Style = $ws.Range("B2").Style
$ws.Range("E10").Value = "  +3.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.02"
$ws.Range("D11").Style = $ws.Range("B2").Style
$ws.Range("E11").Value = "  +2.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0929"
$ws.Range("D12").Style = $ws.Range("B2").Style
$ws.Range("E12").Value = "  -1.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.07"
$ws.Range("D13").Style = $ws.Range("B2").Style
$ws.Range("E13").Value = "  +8.08%  "
$ws.Range("E14").Value = "  +0.69%  "
$ws.Range("D15").Value = "2.500.96"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.868"
$ws.Range("D16").Style = $ws.Range("B2").Style
$ws.Range("E16").Value = "  +5.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.39"
$ws.Range("D17").Style = $ws.Range("B2").Style
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("D18").Value = "2.179.59"
$ws.Range("E18").Value = "  -0.38%  "
$ws.Range("D19").Value = "41.107.76"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("E20").Value = "  +1.30%  "
$ws.Range("E21").Value = "  +1.95%  "
$ws.Range("E22").Value = "  -0.74%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "231.29"
$ws.Range("D23").Style = $ws.Range("B2").Style
$ws.Range("E23").Value = "  +1.18%  "
$ws.Range("E24").Value = "  +1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.94"
$ws.Range("D25").Style = $ws.Range("B2").Style
$ws.Range("E25").Value = "  +10.17%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.69"
$ws.Range("D26").Style = $ws.Range("B2").Style
$ws.Range("E26").Value = "  +21.60%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +6.14%  "
$ws.Range("E29").Value = "  +0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "167.81"
$ws.Range("D30").Style = $ws.Range("B2").Style
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("E31").Value = "  +2.36%  "
$ws.Range("E32").Value = "  -0.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0745"
$ws.Range("D33").Style = $ws.Range("B2").Style
$ws.Range("E33").Value = "  +7.22%  "
$ws.Range("E34").Value = "  +0.25%  "
$ws.Range("E35").Value = "  +6.57%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.27"
$ws.Range("D36").Style = $ws.Range("B2").Style
$ws.Range("E36").Value = "  +12.72%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.14"
$ws.Range("D37").Style = $ws.Range("B2").Style
$ws.Range("E37").Value = "  +8.65%  "
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.59"
$ws.Range("D38").Style = $ws.Range("B2").Style
$ws.Range("E38").Value = "  +1.06%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0298"
$ws.Range("D39").Style = $ws.Range("B2").Style
$ws.Range("E39").Value = "  +11.91%  "
$ws.Range("E40").Value = "  -2.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.43"
$ws.Range("D41").Style = $ws.Range("B2").Style
$ws.Range("E41").Value = "  +21.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.65"
$ws.Range("D42").Style = $ws.Range("B2").Style
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.18"
$ws.Range("D43").Style = $ws.Range("B2").Style
$ws.Range("E43").Value = "  +3.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.03"
$ws.Range("D44").Style = $ws.Range("B2").Style
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("E45").Value = "  +5.76%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.62"
$ws.Range("D46").Style = $ws.Range("B2").Style
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("E47").Value = "  +2.97%  "
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +4.60%  "
$ws.Range("E50").Value = "  +1.23%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.25"
$ws.Range("D51").Style = $ws.Range("B2").Style
$ws.Range("E51").Value = "  -4.15%  "
